$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving the cells
# original style/number-format, forcing Excel to store it as text
# instead of re-interpreting numeric-looking strings (e.g. "1.001")
# as actual numbers.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '30.214.65'
Set-TextValue $ws.Range("E2") '  -0.25%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.863.75'

# Row 4
Set-TextValue $ws.Range("E4") '  -0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '241.74'
Set-TextValue $ws.Range("E5") '  +2.91%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.001'
Set-TextValue $ws.Range("E6") '  +0.02%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.4704'
Set-TextValue $ws.Range("E7") '  +0.17%  '

# Row 8
Set-TextValue $ws.Range("D8") '42.73'
Set-TextValue $ws.Range("E8") '  -0.43%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.2852'
Set-TextValue $ws.Range("E9") '  -0.56%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.06469'
Set-TextValue $ws.Range("E10") '  -2.00%  '

# Row 11
Set-TextValue $ws.Range("D11") '20.78'
Set-TextValue $ws.Range("E11") '  -3.75%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.07698'
Set-TextValue $ws.Range("E12") '  -3.26%  '

# Row 13
Set-TextValue $ws.Range("D13") '1.863.77'
Set-TextValue $ws.Range("E13") '  -0.38%  '

# Row 14
Set-TextValue $ws.Range("D14") '94.93'
Set-TextValue $ws.Range("E14") '  -1.82%  '

# Row 15
Set-TextValue $ws.Range("B15") 'Polkadot'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D15") '5.085'
Set-TextValue $ws.Range("E15") '  -0.44%  '

# Row 16
Set-TextValue $ws.Range("B16") 'Polygon'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D16") '0.6834'
Set-TextValue $ws.Range("E16") '  -2.26%  '

# Row 17
Set-TextValue $ws.Range("D17") '268.92'
Set-TextValue $ws.Range("E17") '  +0.06%  '

# Row 18
Set-TextValue $ws.Range("D18") '30.204.31'
Set-TextValue $ws.Range("E18") '  -0.48%  '

# Row 19
Set-TextValue $ws.Range("E19") '  -5.39%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.000007528'
Set-TextValue $ws.Range("E20") '  -3.42%  '

# Row 21
Set-TextValue $ws.Range("E21") '  -0.01%  '

# Row 22
Set-TextValue $ws.Range("D22") '2.107.03'
Set-TextValue $ws.Range("E22") '  -0.67%  '

# Row 23
Set-TextValue $ws.Range("E23") '  -0.10%  '

# Row 24
Set-TextValue $ws.Range("D24") '5.202'
Set-TextValue $ws.Range("E24") '  -1.30%  '

# Row 25
Set-TextValue $ws.Range("D25") '6.117'
Set-TextValue $ws.Range("E25") '  -1.55%  '

# Row 26
Set-TextValue $ws.Range("D26") '9.330'
Set-TextValue $ws.Range("E26") '  -0.47%  '

# Row 27
Set-TextValue $ws.Range("D27") '165.67'
Set-TextValue $ws.Range("E27") '  -1.01%  '

# Row 28
Set-TextValue $ws.Range("D28") '18.83'
Set-TextValue $ws.Range("E28") '  -0.19%  '

# Row 29
Set-TextValue $ws.Range("E29") '  -2.94%  '

# Row 30
Set-TextValue $ws.Range("E30") '  +0.60%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.09806'
Set-TextValue $ws.Range("E31") '  -1.07%  '

# Row 32
Set-TextValue $ws.Range("D32") '1.504'
Set-TextValue $ws.Range("E32") '  +2.98%  '

# Row 33
Set-TextValue $ws.Range("E33") '  -2.16%  '

# Row 34
Set-TextValue $ws.Range("D34") '3.987'
Set-TextValue $ws.Range("E34") '  -1.69%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.04696'
Set-TextValue $ws.Range("E35") '  -0.33%  '

# Row 36
Set-TextValue $ws.Range("E36") '  -2.03%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.6846'
Set-TextValue $ws.Range("E37") '  -2.55%  '

# Row 38
Set-TextValue $ws.Range("E38") '  -0.66%  '

# Row 39
Set-TextValue $ws.Range("E39") '  -1.64%  '

# Row 40
Set-TextValue $ws.Range("E40") '  -2.19%  '

# Row 41
Set-TextValue $ws.Range("D41") '6.345'
Set-TextValue $ws.Range("E41") '  +1.39%  '

# Row 42
Set-TextValue $ws.Range("D42") '70.40'
Set-TextValue $ws.Range("E42") '  -2.20%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.9998'
Set-TextValue $ws.Range("E43") '  +0.04%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.8391'
Set-TextValue $ws.Range("E44") '  -0.26%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.887'
Set-TextValue $ws.Range("E45") '  -3.65%  '

# Row 46
Set-TextValue $ws.Range("D46") '101.84'
Set-TextValue $ws.Range("E46") '  -0.92%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.4067'
Set-TextValue $ws.Range("E47") '  -2.68%  '

# Row 48
Set-TextValue $ws.Range("D48") '9.146'
Set-TextValue $ws.Range("E48") '  -0.30%  '

# Row 49
Set-TextValue $ws.Range("D49") '6.950'
Set-TextValue $ws.Range("E49") '  -2.21%  '

# Row 50
Set-TextValue $ws.Range("D50") '922.32'
Set-TextValue $ws.Range("E50") '  +0.47%  '

# Row 51
Set-TextValue $ws.Range("D51") '34.45'
Set-TextValue $ws.Range("E51") '  -0.30%  '
